# Applies the textual corrections described in the commit diff:
#   1. "tousiours"                -> "tousjours"          (i -> j)
#   2. "l'escuelle"                -> "lescuelle"           (l' -> l)
#   3-7. "et mesmement au commancement que tu verseras a "
#        (run boundaries reshuffled by the upstream regeneration tool,
#         overall visible text unchanged)
#   8. "j'ay"                      -> "jay"                 (j' -> j)
#
# Strategy: find a short, document-unique "anchor" string that sits
# immediately in front of the run we need to edit. Collapse/advance a
# Range to right after that anchor, then -- instead of trusting raw
# character-offset arithmetic (paragraph marks etc. add invisible units
# to Range start/end positions) -- do a second Find for the exact old
# run text starting at that position and ending at the document end.
# Because the old text is supposed to sit immediately (modulo an
# occasional paragraph-mark) after the anchor, the first hit found is
# the right one; we sanity check the gap is tiny. The replace itself is
# then executed with Find/Replace scoped to exactly that narrow Range,
# so only the single target run is ever touched.

$d = $word.ActiveDocument

function Replace-AfterAnchor {
    param(
        [string]$Anchor,
        [string]$OldText,
        [string]$NewText
    )

    $anchorRange = $d.Content
    $found = $anchorRange.Find.Execute($Anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "ANCHOR NOT FOUND: '$Anchor'"
        return
    }
    $anchorEnd = $anchorRange.End

    $search = $d.Range($anchorEnd, $d.Content.End)
    $found2 = $search.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found2) {
        Write-Output "TARGET TEXT NOT FOUND after anchor '$Anchor': '$OldText'"
        return
    }

    $gap = $search.Start - $anchorEnd
    if ($gap -lt 0 -or $gap -gt 3) {
        Write-Output "TARGET TOO FAR from anchor '$Anchor' (gap=$gap) - aborting to avoid wrong edit"
        return
    }

    $target = $d.Range($search.Start, $search.End)
    $target.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2) | Out-Null
    Write-Output "Replaced '$OldText' -> '$NewText' after anchor '$Anchor' (gap=$gap)"
}

# Processed in reverse document order so that earlier anchors (further to
# the left in the text) are never disturbed by an edit that has not
# happened yet.

Replace-AfterAnchor "e dont " "j’" "j"
Replace-AfterAnchor "t au commanc<exp>emen</exp>" "t que tu verseras a " " que tu verseras a "
Replace-AfterAnchor "t au commanc<exp>" "emen" "ent"
Replace-AfterAnchor "et mesme<exp>men</exp>" "t au commanc" " au commancem"
Replace-AfterAnchor "et mesme<exp>" "men" "ent"
Replace-AfterAnchor "espes<lb/>" "et mesme" "et mesmem"
Replace-AfterAnchor "et co<exp>mm</exp>e " "l’" "l"
Replace-AfterAnchor "haulse tous" "i" "j"
